{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Applies 3 content edits described by the diff:\n//   1. \"Utkarsh Behre, Nikhil Kamath\" -> \"Utkarsh Behre, Nikhil D Kamath\"\n//      (split into 3 runs: \", Nikhil\" + \" D\" + \" Kamath\", same bold formatting)\n//   2. Remove the stray \"_GoBack\" bookmark left after a picture.\n//   3. Clean up proofing (spell-check) marks around \"Acrobot\" by merging the\n//      surrounding runs back into a single run with the same visible text.\n\nconst body = context.document.body;\n\n// --- 1. \"..., Nikhil Kamath\" -> \"..., Nikhil D Kamath\" (3 runs) ---------\nconst nameResults = body.search(\", Nikhil Kamath\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  const nameRange = nameResults.items[0];\n  const nameOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r w:rsidR=\"00A745DD\"><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t>, Nikhil</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t xml:space=\"preserve\"> D</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t xml:space=\"preserve\"> Kamath</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  nameRange.insertOoxml(nameOoxml, \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Remove the \"_GoBack\" bookmark -----------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. Merge the runs/proofing marks around \"Acrobot\" ------------------\nconst acrobotText =\n  \" and also Acrobot environments. We did this for all 5 algorithms to see \" +\n  \"the comparisons. We have a notebook for each environment where we ran \" +\n  \"all algorithms and plotted graphs for them. The exception being PPO \" +\n  \"for which we created a separate notebook since it is a new algorithm we tried.\";\n\nconst acrobotResults = body.search(acrobotText, { matchCase: true });\nacrobotResults.load(\"items\");\nawait context.sync();\n\nif (acrobotResults.items.length > 0) {\n  const acrobotRange = acrobotResults.items[0];\n  const acrobotOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">' + acrobotText + '</w:t></w:r></w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  acrobotRange.insertOoxml(acrobotOoxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Applies 3 content edits described by the diff:\n#   1. \"Utkarsh Behre, Nikhil Kamath\" -> \"Utkarsh Behre, Nikhil D Kamath\"\n#      (split into 3 runs: \", Nikhil\" + \" D\" + \" Kamath\", same bold formatting)\n#   2. Remove the stray \"_GoBack\" bookmark left after a picture.\n#   3. Clean up proofing (spell-check) marks around \"Acrobot\" by merging the\n#      surrounding runs back into a single run with the same visible text.\n\n$d = $word.ActiveDocument\n\n# --- 1. \"..., Nikhil Kamath\" -> \"..., Nikhil D Kamath\" (3 runs) ---------\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Text = \", Nikhil Kamath\"\n$found1 = $find1.Find.Execute()\n\nif ($found1) {\n    $nameRange = $d.Range($find1.Start, $find1.End)\n    $nameXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n        '<w:r w:rsidR=\"00A745DD\"><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t>, Nikhil</w:t></w:r>' +\n        '<w:r><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t xml:space=\"preserve\"> D</w:t></w:r>' +\n        '<w:r><w:rPr><w:b/><w:bCs/><w:spacing w:val=\"5\"/><w:kern w:val=\"1\"/></w:rPr><w:t xml:space=\"preserve\"> Kamath</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n    $nameRange.InsertXML($nameXml)\n}\n\n# --- 2. Remove the \"_GoBack\" bookmark -----------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3. Merge the runs/proofing marks around \"Acrobot\" ------------------\n$acrobotText = \" and also Acrobot environments. We did this for all 5 algorithms to see \" +\n    \"the comparisons. We have a notebook for each environment where we ran \" +\n    \"all algorithms and plotted graphs for them. The exception being PPO \" +\n    \"for which we created a separate notebook since it is a new algorithm we tried.\"\n\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Text = $acrobotText\n$found2 = $find2.Find.Execute()\n\nif ($found2) {\n    $acrobotRange = $d.Range($find2.Start, $find2.End)\n    $acrobotXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p><w:r><w:t xml:space=\"preserve\">' + $acrobotText + '</w:t></w:r></w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n    $acrobotRange.InsertXML($acrobotXml)\n}\n"}
